# minor changes in chapter 7
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update competitor comparison values (column B = GitFit, C = eGym, D = myClub, E = technogym)
$ws.Range("B2").Value = 6
$ws.Range("B3").Value = 7
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 3
$ws.Range("B4").Value = 6
$ws.Range("B5").Value = 6
$ws.Range("B6").Value = 3
$ws.Range("B9").Value = 4

# Update the active selection to match the saved view state
$ws.Activate()
$ws.Range("B6").Select()
